# February 24, 2021 update
# Appends 13 new daily rows (336-348, covering 2021-02-12 .. 2021-02-24)
# to the GermanyAllNumbers worksheet, following the exact same pattern as
# every preceding row: A/B/C are the constant iso_code/continent/location
# strings, D is "previous date + 1", E/G/I/J are rolling-window formulas,
# F and H are the new day's raw "new_cases" / "total_deaths" readings, and
# K-P are all 0 (no data reported for those columns in this period).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GermanyAllNumbers")

# New daily inputs taken straight from the source data update: F = new_cases,
# H = total_deaths for that date. Everything else on the row is derived.
$newRows = @(
    @{ Row = 336; F = 9860;  H = 64191 },
    @{ Row = 337; F = 8354;  H = 64672 },
    @{ Row = 338; F = 6114;  H = 64960 },
    @{ Row = 339; F = 4426;  H = 65076 },
    @{ Row = 340; F = 3856;  H = 65604 },
    @{ Row = 341; F = 7556;  H = 66164 },
    @{ Row = 342; F = 10207; H = 66698 },
    @{ Row = 343; F = 9113;  H = 67206 },
    @{ Row = 344; F = 9164;  H = 67696 },
    @{ Row = 345; F = 7676;  H = 67841 },
    @{ Row = 346; F = 4369;  H = 67903 },
    @{ Row = 347; F = 3883;  H = 68318 },
    @{ Row = 348; F = 8007;  H = 68740 }
)

foreach ($d in $newRows) {
    $r = $d.Row
    $prev = $r - 1
    $winStart = $r - 6

    $ws.Range("A$r").Value = "DEU"
    $ws.Range("B$r").Value = "Europe"
    $ws.Range("C$r").Value = "Germany"

    $ws.Range("D$r").Formula = "=D$prev+1"
    $ws.Range("E$r").Formula = "=E$prev+F$r"
    $ws.Range("F$r").Value = $d.F
    $ws.Range("G$r").Formula = "=SUM(F$winStart`:F$r)/7"
    $ws.Range("H$r").Value = $d.H
    $ws.Range("I$r").Formula = "=H$r-H$prev"
    $ws.Range("J$r").Formula = "=SUM(I$winStart`:I$r)/7"

    $ws.Range("K$r").Value = 0
    $ws.Range("L$r").Value = 0
    $ws.Range("M$r").Value = 0
    $ws.Range("N$r").Value = 0
    $ws.Range("O$r").Value = 0
    $ws.Range("P$r").Value = 0
}

# Recalculate so every new formula cell carries a fresh cached value.
$excel.Calculate()

# Match the author's final cursor/selection position from the diff
# (topLeftCell B325, active cell F342) as closely as the view model allows.
$ws.Range("F342").Select()
